$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# --- Row 15 ---
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = -40

# --- Row 16 ---
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -25
$ws.Range("I16").Value = 48
$ws.Range("J16").Value = 62
$ws.Range("K16").Value = -22.580645161290
$ws.Range("L16").Value = -29.411764705882
$ws.Range("M16").Value = -39.240506329113

# --- Row 17 ---
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 6.25
$ws.Range("I17").Value = 92
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = -8
$ws.Range("L17").Value = -14.814814814814
$ws.Range("M17").Value = 43.75

# --- Row 18 ---
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = -30
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 41
$ws.Range("K18").Value = 29.268292682926
$ws.Range("L18").Value = -18.461538461538
$ws.Range("M18").Value = 47.222222222222

# --- Row 19 ---
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 4
$ws.Range("F19").Value = 16
$ws.Range("H19").Value = -33.333333333333
$ws.Range("I19").Value = 102
$ws.Range("J19").Value = 116
$ws.Range("K19").Value = -12.068965517241
$ws.Range("L19").Value = -8.928571428571
$ws.Range("M19").Value = 5.154639175257

# --- Row 20 ---
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 8
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 45
$ws.Range("J20").Value = 52
$ws.Range("K20").Value = -13.461538461538
$ws.Range("L20").Value = 150
$ws.Range("M20").Value = 66.666666666666

# --- Row 21 ---
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 17
$ws.Range("E21").Value = -17.647058823529
$ws.Range("F21").Value = 54
$ws.Range("G21").Value = 68
$ws.Range("H21").Value = -20.588235294117
$ws.Range("I21").Value = 344
$ws.Range("J21").Value = 377
$ws.Range("K21").Value = -8.753315649867
$ws.Range("L21").Value = -8.021390374331
$ws.Range("M21").Value = 10.610932475884

# --- Row 22 ---
# D22/E22 convert from the "no data" placeholder text to real numbers, so the
# number format has to be (re)applied to match the numeric sibling columns.
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = -100
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I22").Value = 4
$ws.Range("J22").Value = 13
$ws.Range("K22").Value = -69.230769230769
$ws.Range("L22").Value = -20
$ws.Range("M22").Value = 33.333333333333

# --- Row 23 ---
$ws.Range("G23").Value = 1

# --- Row 24 ---
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 29
$ws.Range("E24").Value = -41.379310344827
$ws.Range("F24").Value = 64
$ws.Range("G24").Value = 109
$ws.Range("H24").Value = -41.284403669724
$ws.Range("I24").Value = 346
$ws.Range("J24").Value = 589
$ws.Range("K24").Value = -41.256366723259
$ws.Range("L24").Value = -19.347319347319
$ws.Range("M24").Value = 92.222222222222

# --- Row 25 ---
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 266.666666666667
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 150
$ws.Range("J25").Value = 136
$ws.Range("K25").Value = 10.294117647058
$ws.Range("L25").Value = 17.1875
$ws.Range("M25").Value = -1.960784313725

# --- Row 26 ---
$ws.Range("G26").Value = 3
$ws.Range("J26").Value = 8
$ws.Range("K26").Value = -37.5

# --- Row 27 ---
# C27 converts from the "no data" placeholder text to a real number.
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -37.5
$ws.Range("I27").Value = 16
$ws.Range("J27").Value = 21
$ws.Range("K27").Value = -23.809523809523
$ws.Range("L27").Value = 33.333333333333

# --- Row 30 ---
# D30/E30 and G30/H30 convert from the "no data" placeholder text to real numbers.
$ws.Range("D30").Value = 1
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").Value = -100
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("G30").Value = 1
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("H30").Value = -100
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("J30").Value = 2
$ws.Range("K30").Value = 0
